$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before column D ("mass.init" shifts right,
#    making room for the new "interval" column).
$ws.Columns("D").Insert()

# 2. New header for the inserted column.
$ws.Range("D4").Value = "interval"

# 3. Row 5 (first / baseline measurement) -------------------------------
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = 385.27
$ws.Range("F5").Value = 385.27
$ws.Range("G5").Formula = "=E5-F5"
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()

# 4. Row 6 ----------------------------------------------------------------
$ws.Range("C6").Value = 0.90690000000000004
$ws.Range("D6").Formula = "=C6-C5"
$ws.Range("E6").Value = 385.26
$ws.Range("F6").Value = 385.04
$ws.Range("G6").Formula = "=E6-F6"
$ws.Range("H6").Formula = "=F5-E6"
$ws.Range("I6").Formula = "=F5-F6"
$ws.Range("J6").ClearContents()

# 5. Row 7 ----------------------------------------------------------------
$ws.Range("C7").Value = 2.0381
$ws.Range("D7").Formula = "=C7-C6"
$ws.Range("E7").Value = 385.04
$ws.Range("F7").Value = 384.8
$ws.Range("G7").Formula = "=E7-F7"
$ws.Range("H7").Formula = "=F6-E7"
$ws.Range("I7").Formula = "=F6-F7"
$ws.Range("J7").ClearContents()

# 6. Column width for the newly inserted "interval" column (D), matching
#    the width used for the "elapsed.time"/"id" style columns.
$ws.Columns("D").ColumnWidth = 12.7109375

# 7. View tweaks captured in the diff (scrolled right, selection moved).
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("I12").Select()

$wb.Save()
